$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add prolificid in rank to use in binary"
# The ranking table (Sheet1, rows 2-25) is regenerated from a refreshed
# realeffort/prolific-id draw: the realeffort score (col F / re_rank) is
# recomputed for every worker, which reshuffles a few same-gender ties
# (their worker-id "index" in col C, name in col D, and race in col G move
# together to the row matching their new re_rank position), while the
# level_0/index-within-gender/gender/rank columns (A, B, E, H) stay fixed.

$index      = @(41, 19, 2, 3, 34, 44, 22, 35, 33, 21, 32, 30, 44, 3, 27, 30, 22, 32, 26, 2, 33, 49, 29, 50)
$names      = @("Giana", "Jewel", "Colleen", "Annes", "Tina", "Nansi", "Khushi", "Lori", "Shaniek", "Bri", "Kellie", "Shadaisia", "Myles", "Quinterius", "Drew", "Matthew", "Edosagbe", "Jamarii", "Juan", "Corey", "Brennan", "Masuf", "Eli", "Damian")
$realeffort = @(7.240540192629654, 6.378978103426058, 6.143455313863114, 5.419772607443591, 5.312796240675778, 4.071991992584385, 1.233832614214271, 1.002782814522061, 0.4746561773749075, 0.4664235049697223, 0.2827895313987391, 0.27386664857579, 13.45172621041747, 8.467383315037575, 7.263320786645187, 7.045923228846132, 5.285624560074965, 5.22667163757618, 5.193444245373518, 4.163691280357252, 4.093901744365527, 3.064644559899139, 2.350791450174602, 2.331889986248744)
$race       = @("White", "Black or African American", "White", "Asian", "White", "Asian", "Asian", "White", "Black or African American", "Black or African American", "White", "Black or African American", "Black or African American", "Black or African American", "White", "White", "Black or African American", "Black or African American", "Hispanic", "White", "White", "Asian", "White", "Black or African American")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $index[$i]
    $ws.Cells.Item($row, 4).Value = $names[$i]
    $ws.Cells.Item($row, 6).Value = $realeffort[$i]
    $ws.Cells.Item($row, 7).Value = $race[$i]
}

Write-Host "Updated realeffort/re_rank table for rows 2-25"
